$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pravin, Shaiju and Tintu paid their dues (400 column contribution of 300 in F column)
$ws.Range("F16").Value = 300
$ws.Range("F23").Value = 300
$ws.Range("F25").Value = 300

# Update the active selection / scroll position to reflect the latest entry (F23)
$ws.Range("F23").Select()
$excel.ActiveWindow.ScrollRow = 1
